$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("Анакина Надежда"): fill in grades.
# C6:E6 already carry the right "data cell" style, just need values.
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 5

# F6:H6 are new entries that need the same border formatting already used
# for the extra-grade columns elsewhere on the sheet (e.g. F23), so copy
# that formatting in before writing the values.
$ws.Range("F23").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F23").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("F23").Copy()
$ws.Range("H6").PasteSpecial(-4122)

$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 5

# Row 17 ("Муллаянова Карина"): fill in grades (style already in place).
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 5

# Move the active selection to F17, matching the saved workbook's view state.
$ws.Range("F17").Select()
